$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-7
# from serial date 45183 (2023-09-14) to 45184 (2023-09-15)
$ws.Range("C2:C7").Value = 45184
